$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift the existing 20 data rows (A2:C21) down by 5 rows (to A7:C26) ---
# Capture the current values first, then write them to the new location. This
# avoids Excel's Range.Insert() formatting inheritance quirks (it would copy
# the header row's style onto the newly inserted blank rows).
$existing = $ws.Range("A2:C21").Value2
$ws.Range("A7:C26").Value2 = $existing

# --- Step 2: write the 5 new rows at the top (A2:C6) ---
$topRows = New-Object 'object[,]' 5,3
$topRows[0,0] = -0.09269879758358;     $topRows[0,1] = 0.0145080499351024;   $topRows[0,2] = -0.007177666760981
$topRows[1,0] = -0.090408056974411;    $topRows[1,1] = -0.0394008085131645;  $topRows[1,2] = -0.0087048299610614
$topRows[2,0] = -0.07376197725534429;  $topRows[2,1] = -0.102472648024559;   $topRows[2,2] = 0.0461203269660472
$topRows[3,0] = -0.1244637966156005;   $topRows[3,1] = -0.4952589869499206;  $topRows[3,2] = 0.25641068816185
$topRows[4,0] = 0.955545961856842;     $topRows[4,1] = 0.418595403432846;    $topRows[4,2] = 0.5012149214744568
$ws.Range("A2:C6").Value2 = $topRows

# --- Step 3: append the 5 new rows at the bottom (A27:C31) ---
$bottomRows = New-Object 'object[,]' 5,3
$bottomRows[0,0] = 0.4641048610210418;  $bottomRows[0,1] = 0.3608686327934265;  $bottomRows[0,2] = 0.3602577745914459
$bottomRows[1,0] = 0.2412917762994766;  $bottomRows[1,1] = 0.2144137024879455;  $bottomRows[1,2] = -0.0186313893646001
$bottomRows[2,0] = 0.1291979998350143;  $bottomRows[2,1] = 0.1751656085252761;  $bottomRows[2,2] = -0.070249505341053
$bottomRows[3,0] = -0.1421788930892944; $bottomRows[3,1] = -0.0774271711707115; $bottomRows[3,2] = -0.0045814891345798
$bottomRows[4,0] = -0.067195177078247;  $bottomRows[4,1] = -0.219300627708435;  $bottomRows[4,2] = -0.2028072625398636
$ws.Range("A27:C31").Value2 = $bottomRows
